# "add delphi navigation" — insert a new "delphi" worksheet (before "ithaca")
# populated with navigation test data, and clear the previously-active
# "specialist" sheet's selection/tab state in favour of the new sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new "delphi" sheet immediately before "ithaca" so sheet
#    order becomes: male, female, specialist, delphi, ithaca.
# ---------------------------------------------------------------------
$ithaca = $wb.Worksheets.Item("ithaca")
$delphi = $wb.Worksheets.Add($ithaca)
$delphi.Name = "delphi"

# ---------------------------------------------------------------------
# 2. Populate header row (row 1) and data row (row 2).
# ---------------------------------------------------------------------
$headers = @("testCaseName","HeaderName","DashboardField","DashboardLink","LinkName1","LinkName2","LinkName3","LinkName4","LinkName5","ApplicationsearchText","LocationText","DataText","SurveyModuleText","QuestionText")
$data    = @("delphiNavigation","Site Administration","Dashboard","Dashboard","Applications","Location themes","Questions","Surveys","Survey modules","Family","location test","Location data test","Survey","Survey")

for ($i = 0; $i -lt $headers.Count; $i++) {
    $delphi.Cells.Item(1, $i + 1).Value = $headers[$i]
    $delphi.Cells.Item(2, $i + 1).Value = $data[$i]
}

# Header row formatting: reuse the workbook's existing header style
# (Calibri 12, black) and a slightly taller row height.
$headerRange = $delphi.Range("A1:N1")
$headerRange.Font.Name = "Calibri"
$headerRange.Font.Size = 12
$delphi.Rows.Item(1).RowHeight = 15.6

# Column widths for A:D.
$delphi.Columns.Item(1).ColumnWidth = 14.5
$delphi.Columns.Item(2).ColumnWidth = 14.833333333333334
$delphi.Columns.Item(3).ColumnWidth = 16
$delphi.Columns.Item(4).ColumnWidth = 14.5

# ---------------------------------------------------------------------
# 3. "specialist" is no longer the active tab — move its selection back
#    to a plain A1:D1 range (matches the committed view state).
# ---------------------------------------------------------------------
$specialist = $wb.Worksheets.Item("specialist")
$specialist.Range("A1:D1").Select() | Out-Null

# ---------------------------------------------------------------------
# 4. "delphi" becomes the active sheet/tab with P11 selected.
# ---------------------------------------------------------------------
$delphi.Activate() | Out-Null
$delphi.Range("P11").Select() | Out-Null
